# Adds the "Implementations" / "User Stories" / "Demo Video" work items
# (rows 8-12) to the Work Breakdown Structure sheet, matching the upstream
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write cell text values --------------------------------------------
# Order matters: new (not-yet-seen) strings are appended to the shared
# string table in the order they are first written, and we want that order
# to match the upstream workbook's sharedStrings.xml (indices 26-34).
$ws.Range("A9").Value  = "User Stories"
$ws.Range("A12").Value = "Demo Video"
$ws.Range("B8").Value  = "Elaborate 1 use-case diagram for each implementation"
$ws.Range("B9").Value  = "Elaborate user stories for each implementation"
$ws.Range("C9").Value  = "Elaborate acceptance criteria for each user story"
$ws.Range("B10").Value = "Study base JabRef code"
$ws.Range("C10").Value = "Locate important classes to modify"
$ws.Range("A10").Value = "Implementations"
$ws.Range("A11").Value = "Unit Tests"

# Cells that reuse strings already present in the workbook.
$ws.Range("A8").Value  = "Use-Case Diagram"
$ws.Range("C8").Value  = "Elaborate a use-case description file: name,description and actors."
$ws.Range("D8").Value  = "Push to remote"
$ws.Range("D9").Value  = "Push to remote"

# --- 2. Clone formatting from existing, visually-equivalent cells ---------
# Each paste is done on a single cell at a time (multi-area / block pastes
# behave inconsistently), and $excel.CutCopyMode is cleared after every
# copy so the marching-ants selection doesn't leak into the next operation.

# Column A category cells (bold + underline, 14pt) - copy from A4.
foreach ($addr in @("A8", "A9", "A10", "A11", "A12")) {
    $ws.Range("A4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Task cells that use the D4/D5 style (s=7) - copy from D4.
foreach ($addr in @("B8", "B9", "B10", "C9", "C10")) {
    $ws.Range("D4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# "C" task description cell that shares the B/E-column look (s=6) - copy
# from B4.
$ws.Range("B4").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# "Push to remote" styled cells (s=10) - copy from F4. E8 stays blank but
# keeps the same style applied.
foreach ($addr in @("D8", "D9", "E8")) {
    $ws.Range("F4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# --- 3. New red-text styles for the still-empty D10/E10 cells -------------
# D10 mirrors the D4-style xf (s=7) but with red font; E10 mirrors the
# B4-style xf (s=6) but with red font.
$ws.Range("D4").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D10").Font.Name = "Calibri"
$ws.Range("D10").Font.Size = 14
$ws.Range("D10").Font.Color = 255

$ws.Range("B4").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E10").Font.Name = "Calibri"
$ws.Range("E10").Font.Size = 14
$ws.Range("E10").Font.Color = 255

# --- 4. Row heights for the newly-added rows -------------------------------
$ws.Rows.Item(8).RowHeight  = 33.6
$ws.Rows.Item(9).RowHeight  = 26.4
$ws.Rows.Item(10).RowHeight = 32.4
$ws.Rows.Item(11).RowHeight = 28.8
$ws.Rows.Item(12).RowHeight = 32.4

# --- 5. Selection matches the saved view in the upstream file -------------
$ws.Range("B11").Select()
